$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-12-02 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-12-03 Sunday", 2) | Out-Null
$d.Content.Find.Execute("38-28=", $true, $false, $false, $false, $false, $true, 1, $false, "53+35=", 2) | Out-Null
$d.Content.Find.Execute("18-8=", $true, $false, $false, $false, $false, $true, 1, $false, "9+3=", 2) | Out-Null
$d.Content.Find.Execute("42+3=", $true, $false, $false, $false, $false, $true, 1, $false, "25-3=", 2) | Out-Null
$d.Content.Find.Execute("74-46=", $true, $false, $false, $false, $false, $true, 1, $false, "34-24=", 2) | Out-Null
$d.Content.Find.Execute("78+17=", $true, $false, $false, $false, $false, $true, 1, $false, "14-1=", 2) | Out-Null
$d.Content.Find.Execute("63+12=", $true, $false, $false, $false, $false, $true, 1, $false, "26+12=", 2) | Out-Null
$d.Content.Find.Execute("24+61=", $true, $false, $false, $false, $false, $true, 1, $false, "59+12=", 2) | Out-Null
$d.Content.Find.Execute("37+13=", $true, $false, $false, $false, $false, $true, 1, $false, "18+6=", 2) | Out-Null
$d.Content.Find.Execute("43-31=", $true, $false, $false, $false, $false, $true, 1, $false, "95-48=", 2) | Out-Null
$d.Content.Find.Execute("29+6=", $true, $false, $false, $false, $false, $true, 1, $false, "7+90=", 2) | Out-Null
$d.Content.Find.Execute("29+21=", $true, $false, $false, $false, $false, $true, 1, $false, "73-52=", 2) | Out-Null
$d.Content.Find.Execute("56+34=", $true, $false, $false, $false, $false, $true, 1, $false, "43-38=", 2) | Out-Null
$d.Content.Find.Execute("23+28=", $true, $false, $false, $false, $false, $true, 1, $false, "0+44=", 2) | Out-Null
$d.Content.Find.Execute("96-12=", $true, $false, $false, $false, $false, $true, 1, $false, "31+1=", 2) | Out-Null
$d.Content.Find.Execute("64-17=", $true, $false, $false, $false, $false, $true, 1, $false, "75+3=", 2) | Out-Null
$d.Content.Find.Execute("84-0=", $true, $false, $false, $false, $false, $true, 1, $false, "37+28=", 2) | Out-Null
$d.Content.Find.Execute("90-22=", $true, $false, $false, $false, $false, $true, 1, $false, "99-34=", 2) | Out-Null
$d.Content.Find.Execute("32-18=", $true, $false, $false, $false, $false, $true, 1, $false, "45+41=", 2) | Out-Null
$d.Content.Find.Execute("46-38=", $true, $false, $false, $false, $false, $true, 1, $false, "45+46=", 2) | Out-Null
$d.Content.Find.Execute("62-49=", $true, $false, $false, $false, $false, $true, 1, $false, "8-1=", 2) | Out-Null
$d.Content.Find.Execute("79-35=", $true, $false, $false, $false, $false, $true, 1, $false, "44-36=", 2) | Out-Null
$d.Content.Find.Execute("79+19=", $true, $false, $false, $false, $false, $true, 1, $false, "63-11=", 2) | Out-Null
$d.Content.Find.Execute("15-10=", $true, $false, $false, $false, $false, $true, 1, $false, "82-39=", 2) | Out-Null
$d.Content.Find.Execute("57+23=", $true, $false, $false, $false, $false, $true, 1, $false, "81-43=", 2) | Out-Null
$d.Content.Find.Execute("66-22=", $true, $false, $false, $false, $false, $true, 1, $false, "92-12=", 2) | Out-Null
$d.Content.Find.Execute("64+33=", $true, $false, $false, $false, $false, $true, 1, $false, "53-51=", 2) | Out-Null
$d.Content.Find.Execute("87-59=", $true, $false, $false, $false, $false, $true, 1, $false, "95-72=", 2) | Out-Null
$d.Content.Find.Execute("26+30=", $true, $false, $false, $false, $false, $true, 1, $false, "20+70=", 2) | Out-Null
$d.Content.Find.Execute("57-2=", $true, $false, $false, $false, $false, $true, 1, $false, "54+12=", 2) | Out-Null
$d.Content.Find.Execute("13+13=", $true, $false, $false, $false, $false, $true, 1, $false, "1+88=", 2) | Out-Null
$d.Content.Find.Execute("90-45=", $true, $false, $false, $false, $false, $true, 1, $false, "91-14=", 2) | Out-Null
$d.Content.Find.Execute("28-22=", $true, $false, $false, $false, $false, $true, 1, $false, "99-15=", 2) | Out-Null
$d.Content.Find.Execute("47-17=", $true, $false, $false, $false, $false, $true, 1, $false, "3+17=", 2) | Out-Null
$d.Content.Find.Execute("91-55=", $true, $false, $false, $false, $false, $true, 1, $false, "72-11=", 2) | Out-Null
$d.Content.Find.Execute("49+27=", $true, $false, $false, $false, $false, $true, 1, $false, "61+4=", 2) | Out-Null
$d.Content.Find.Execute("60-27=", $true, $false, $false, $false, $false, $true, 1, $false, "7+48=", 2) | Out-Null
$d.Content.Find.Execute("58-12=", $true, $false, $false, $false, $false, $true, 1, $false, "22+69=", 2) | Out-Null
$d.Content.Find.Execute("78+7=", $true, $false, $false, $false, $false, $true, 1, $false, "59-50=", 2) | Out-Null
$d.Content.Find.Execute("75+11=", $true, $false, $false, $false, $false, $true, 1, $false, "22+15=", 2) | Out-Null
$d.Content.Find.Execute("46-35=", $true, $false, $false, $false, $false, $true, 1, $false, "27+72=", 2) | Out-Null
$d.Content.Find.Execute("74-42=", $true, $false, $false, $false, $false, $true, 1, $false, "62+9=", 2) | Out-Null
$d.Content.Find.Execute("3+58=", $true, $false, $false, $false, $false, $true, 1, $false, "46-43=", 2) | Out-Null
$d.Content.Find.Execute("51-36=", $true, $false, $false, $false, $false, $true, 1, $false, "48-42=", 2) | Out-Null
$d.Content.Find.Execute("27+63=", $true, $false, $false, $false, $false, $true, 1, $false, "62-47=", 2) | Out-Null
$d.Content.Find.Execute("93-71=", $true, $false, $false, $false, $false, $true, 1, $false, "51-9=", 2) | Out-Null
$d.Content.Find.Execute("71+10=", $true, $false, $false, $false, $false, $true, 1, $false, "64+35=", 2) | Out-Null
$d.Content.Find.Execute("84-80=", $true, $false, $false, $false, $false, $true, 1, $false, "49-37=", 2) | Out-Null
$d.Content.Find.Execute("68-24=", $true, $false, $false, $false, $false, $true, 1, $false, "59-7=", 2) | Out-Null
$d.Content.Find.Execute("90-71=", $true, $false, $false, $false, $false, $true, 1, $false, "80-4=", 2) | Out-Null
$d.Content.Find.Execute("4+56=", $true, $false, $false, $false, $false, $true, 1, $false, "40+3=", 2) | Out-Null
$d.Content.Find.Execute("14+28=", $true, $false, $false, $false, $false, $true, 1, $false, "78+15=", 2) | Out-Null
$d.Content.Find.Execute("83-5=", $true, $false, $false, $false, $false, $true, 1, $false, "41-37=", 2) | Out-Null
$d.Content.Find.Execute("34-32=", $true, $false, $false, $false, $false, $true, 1, $false, "7+7=", 2) | Out-Null
$d.Content.Find.Execute("60-16=", $true, $false, $false, $false, $false, $true, 1, $false, "57+22=", 2) | Out-Null
$d.Content.Find.Execute("51+28=", $true, $false, $false, $false, $false, $true, 1, $false, "90+9=", 2) | Out-Null
$d.Content.Find.Execute("30+66=", $true, $false, $false, $false, $false, $true, 1, $false, "80+4=", 2) | Out-Null
$d.Content.Find.Execute("31+36=", $true, $false, $false, $false, $false, $true, 1, $false, "64-59=", 2) | Out-Null
$d.Content.Find.Execute("27-6=", $true, $false, $false, $false, $false, $true, 1, $false, "42-37=", 2) | Out-Null
$d.Content.Find.Execute("3+26=", $true, $false, $false, $false, $false, $true, 1, $false, "35-31=", 2) | Out-Null
$d.Content.Find.Execute("49+21=", $true, $false, $false, $false, $false, $true, 1, $false, "70-26=", 2) | Out-Null
$d.Content.Find.Execute("73-71=", $true, $false, $false, $false, $false, $true, 1, $false, "44-32=", 2) | Out-Null
$d.Content.Find.Execute("52-40=", $true, $false, $false, $false, $false, $true, 1, $false, "46+40=", 2) | Out-Null
$d.Content.Find.Execute("76-72=", $true, $false, $false, $false, $false, $true, 1, $false, "96-39=", 2) | Out-Null
$d.Content.Find.Execute("52-42=", $true, $false, $false, $false, $false, $true, 1, $false, "61+2=", 2) | Out-Null
$d.Content.Find.Execute("92-80=", $true, $false, $false, $false, $false, $true, 1, $false, "4+82=", 2) | Out-Null
$d.Content.Find.Execute("7+29=", $true, $false, $false, $false, $false, $true, 1, $false, "25-10=", 2) | Out-Null
$d.Content.Find.Execute("17+2=", $true, $false, $false, $false, $false, $true, 1, $false, "29+4=", 2) | Out-Null
$d.Content.Find.Execute("17+61=", $true, $false, $false, $false, $false, $true, 1, $false, "95+0=", 2) | Out-Null
$d.Content.Find.Execute("56-21=", $true, $false, $false, $false, $false, $true, 1, $false, "5+30=", 2) | Out-Null
$d.Content.Find.Execute("52+36=", $true, $false, $false, $false, $false, $true, 1, $false, "60+22=", 2) | Out-Null
$d.Content.Find.Execute("52-11=", $true, $false, $false, $false, $false, $true, 1, $false, "33+19=", 2) | Out-Null
$d.Content.Find.Execute("56-19=", $true, $false, $false, $false, $false, $true, 1, $false, "53+24=", 2) | Out-Null
$d.Content.Find.Execute("84-38=", $true, $false, $false, $false, $false, $true, 1, $false, "86+6=", 2) | Out-Null
$d.Content.Find.Execute("1+45=", $true, $false, $false, $false, $false, $true, 1, $false, "10+49=", 2) | Out-Null
$d.Content.Find.Execute("30+42=", $true, $false, $false, $false, $false, $true, 1, $false, "52-13=", 2) | Out-Null
$d.Content.Find.Execute("95-62=", $true, $false, $false, $false, $false, $true, 1, $false, "76-73=", 2) | Out-Null
$d.Content.Find.Execute("66-14=", $true, $false, $false, $false, $false, $true, 1, $false, "22-18=", 2) | Out-Null
$d.Content.Find.Execute("52+0=", $true, $false, $false, $false, $false, $true, 1, $false, "75-47=", 2) | Out-Null
$d.Content.Find.Execute("60+34=", $true, $false, $false, $false, $false, $true, 1, $false, "73-69=", 2) | Out-Null
$d.Content.Find.Execute("88-4=", $true, $false, $false, $false, $false, $true, 1, $false, "88-63=", 2) | Out-Null
$d.Content.Find.Execute("7+38=", $true, $false, $false, $false, $false, $true, 1, $false, "92-56=", 2) | Out-Null
$d.Content.Find.Execute("47-19=", $true, $false, $false, $false, $false, $true, 1, $false, "93-88=", 2) | Out-Null
$d.Content.Find.Execute("35+2=", $true, $false, $false, $false, $false, $true, 1, $false, "55+4=", 2) | Out-Null
$d.Content.Find.Execute("52-48=", $true, $false, $false, $false, $false, $true, 1, $false, "97-64=", 2) | Out-Null
$d.Content.Find.Execute("74-62=", $true, $false, $false, $false, $false, $true, 1, $false, "52+28=", 2) | Out-Null
$d.Content.Find.Execute("46-22=", $true, $false, $false, $false, $false, $true, 1, $false, "26+5=", 2) | Out-Null
$d.Content.Find.Execute("61-15=", $true, $false, $false, $false, $false, $true, 1, $false, "31+46=", 2) | Out-Null
$d.Content.Find.Execute("25+20=", $true, $false, $false, $false, $false, $true, 1, $false, "85-85=", 2) | Out-Null
$d.Content.Find.Execute("85-79=", $true, $false, $false, $false, $false, $true, 1, $false, "67+7=", 2) | Out-Null
$d.Content.Find.Execute("4+30=", $true, $false, $false, $false, $false, $true, 1, $false, "24+4=", 2) | Out-Null
$d.Content.Find.Execute("91-84=", $true, $false, $false, $false, $false, $true, 1, $false, "30-3=", 2) | Out-Null
$d.Content.Find.Execute("33+59=", $true, $false, $false, $false, $false, $true, 1, $false, "22+72=", 2) | Out-Null
$d.Content.Find.Execute("86-43=", $true, $false, $false, $false, $false, $true, 1, $false, "84-7=", 2) | Out-Null
$d.Content.Find.Execute("0+11=", $true, $false, $false, $false, $false, $true, 1, $false, "4+47=", 2) | Out-Null
$d.Content.Find.Execute("95-44=", $true, $false, $false, $false, $false, $true, 1, $false, "13+83=", 2) | Out-Null
$d.Content.Find.Execute("66-16=", $true, $false, $false, $false, $false, $true, 1, $false, "4+55=", 2) | Out-Null
$d.Content.Find.Execute("79-27=", $true, $false, $false, $false, $false, $true, 1, $false, "93-92=", 2) | Out-Null
$d.Content.Find.Execute("64-57=", $true, $false, $false, $false, $false, $true, 1, $false, "48+49=", 2) | Out-Null
$d.Content.Find.Execute("68+29=", $true, $false, $false, $false, $false, $true, 1, $false, "91-82=", 2) | Out-Null
$d.Content.Find.Execute("70+17=", $true, $false, $false, $false, $false, $true, 1, $false, "45-15=", 2) | Out-Null
